$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the y_nrSteps column (G) where the original value was erroneously "-0";
# normalize these to plain 0.
$zeroRows = @(3, 11, 15, 21, 25, 28)
foreach ($r in $zeroRows) {
    $ws.Cells.Item($r, 7).Value = 0
}

# Update rows where y_corrSteps (E), y_nrSteps (G) and alienID (H) needed correction.
$adjustRows = @(4, 8, 16, 18, 23, 27)
foreach ($r in $adjustRows) {
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 5).Value = $eVal - 1
    $ws.Cells.Item($r, 7).Value = -3
    $ws.Cells.Item($r, 8).Value = 13
}
